$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($addr, $val) {
    $c = $ws.Range($addr)
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Clear stray empty inline-string markers in I418:I420 (collapse to truly empty cells) ---
$ws.Range("I418").Value = ""
$ws.Range("I419").Value = ""
$ws.Range("I420").Value = ""

# --- Append new sprint S34 rows (421-449) ---
# Row 421
Set-CellText A421 'S34'
Set-CellText B421 'G01'
Set-CellText C421 'Holdings Exit Automation MVP: DB schema + models'
Set-CellText D421 'S34_G01_TB001'
Set-CellText E421 'Backend/DB: Add SQLAlchemy models for holding exit subscriptions + events (states, triggers, sizing, audit).'
Set-CellText G421 'planned'

# Row 422
Set-CellText A422 'S34'
Set-CellText B422 'G01'
Set-CellText C422 'Holdings Exit Automation MVP: DB schema + models'
Set-CellText D422 'S34_G01_TB002'
Set-CellText E422 'Backend/DB: Add Alembic migration to create holding_exit_subscriptions and holding_exit_events with indexes + CHECK constraints (cross-dialect safe).'
Set-CellText G422 'planned'

# Row 423
Set-CellText A423 'S34'
Set-CellText B423 'G01'
Set-CellText C423 'Holdings Exit Automation MVP: DB schema + models'
Set-CellText D423 'S34_G01_TB003'
Set-CellText E423 'Backend/DB: Add symbol identity normalization helper for holdings scope (exchange+symbol canonical form) and reuse it in all related features.'
Set-CellText G423 'planned'

# Row 424
Set-CellText A424 'S34'
Set-CellText B424 'G01'
Set-CellText C424 'Holdings Exit Automation MVP: DB schema + models'
Set-CellText D424 'S34_G01_TB004'
Set-CellText E424 'Backend/DB: Define constants/enums for trigger kinds, size modes, statuses; ensure consistent values across Pydantic, DB checks, and UI.'
Set-CellText G424 'planned'

# Row 425
Set-CellText A425 'S34'
Set-CellText B425 'G02'
Set-CellText C425 'Holdings Exit Automation MVP: API + audit + safety posture'
Set-CellText D425 'S34_G02_TB001'
Set-CellText E425 'Backend/API: Implement CRUD endpoints for holding exit subscriptions (list/create/patch/pause/resume/delete).'
Set-CellText G425 'planned'

# Row 426
Set-CellText A426 'S34'
Set-CellText B426 'G02'
Set-CellText C426 'Holdings Exit Automation MVP: API + audit + safety posture'
Set-CellText D426 'S34_G02_TB002'
Set-CellText E426 'Backend/API: Implement events endpoint (list events by subscription) and event writer helper (append-only).'
Set-CellText G426 'planned'

# Row 427
Set-CellText A427 'S34'
Set-CellText B427 'G02'
Set-CellText C427 'Holdings Exit Automation MVP: API + audit + safety posture'
Set-CellText D427 'S34_G02_TB003'
Set-CellText E427 'Backend/API: Validation rules for MVP (CNC-only, SELL-only, percent bounds, qty bounds, trigger kind/value validation, dispatch_mode=MANUAL by default).'
Set-CellText G427 'planned'

# Row 428
Set-CellText A428 'S34'
Set-CellText B428 'G02'
Set-CellText C428 'Holdings Exit Automation MVP: API + audit + safety posture'
Set-CellText D428 'S34_G02_TB004'
Set-CellText E428 'Backend/API: Pause/Resume semantics (resume clears pending_order_id + recomputes next_eval_at; pause records reason).'
Set-CellText G428 'planned'

# Row 429
Set-CellText A429 'S34'
Set-CellText B429 'G02'
Set-CellText C429 'Holdings Exit Automation MVP: API + audit + safety posture'
Set-CellText D429 'S34_G02_TB005'
Set-CellText E429 'Backend/API: Add feature-flag gating (ST_HOLDINGS_EXIT_ENABLED) + optional per-user/per-symbol allowlist hook (safe rollout).'
Set-CellText G429 'planned'

# Row 430
Set-CellText A430 'S34'
Set-CellText B430 'G03'
Set-CellText C430 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D430 'S34_G03_TB001'
Set-CellText E430 'Backend/Engine: Add holdings exit engine scheduler (startup hook + periodic loop) with per-subscription next_eval_at scheduling (adaptive zones).'
Set-CellText G430 'planned'

# Row 431
Set-CellText A431 'S34'
Set-CellText B431 'G03'
Set-CellText C431 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D431 'S34_G03_TB002'
Set-CellText E431 'Backend/Engine: Batch fetch broker holdings by (user, broker) and batch fetch quotes for candidate symbols; handle outages by skipping cycle + emitting EVAL_SKIPPED events.'
Set-CellText G431 'planned'

# Row 432
Set-CellText A432 'S34'
Set-CellText B432 'G03'
Set-CellText C432 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D432 'S34_G03_TB003'
Set-CellText E432 'Backend/Engine: Implement trigger evaluator for TARGET_ABS_PRICE + TARGET_PCT_FROM_AVG_BUY; compute next_eval_at (Far/Near/VeryNear) deterministically.'
Set-CellText G432 'planned'

# Row 433
Set-CellText A433 'S34'
Set-CellText B433 'G03'
Set-CellText C433 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D433 'S34_G03_TB004'
Set-CellText E433 'Backend/Engine: Implement quantity resolver (ABS_QTY, PCT_OF_POSITION) with integer clamping to broker holdings qty (CNC) + min_qty handling.'
Set-CellText G433 'planned'

# Row 434
Set-CellText A434 'S34'
Set-CellText B434 'G03'
Set-CellText C434 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D434 'S34_G03_TB005'
Set-CellText E434 'Backend/Engine: Create WAITING MANUAL SELL order on trigger (is_exit=true, client_order_id prefix HEX:, clear error_message context) and transition subscription to ORDER_CREATED with pending_order_id.'
Set-CellText G434 'planned'

# Row 435
Set-CellText A435 'S34'
Set-CellText B435 'G03'
Set-CellText C435 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D435 'S34_G03_TB006'
Set-CellText E435 'Backend/Engine: Exit arbiter integration (in-flight exit detection via order status set); when conflict, still create WAITING order with annotation + EXIT_QUEUED event (subscription priority via labeling).'
Set-CellText G435 'planned'

# Row 436
Set-CellText A436 'S34'
Set-CellText B436 'G03'
Set-CellText C436 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D436 'S34_G03_TB007'
Set-CellText E436 'Backend/Engine: Idempotency + restart safety (trigger_key, pending_order_id, unique client_order_id) and DB locking strategy (single worker on SQLite; SKIP LOCKED on Postgres).'
Set-CellText G436 'planned'

# Row 437
Set-CellText A437 'S34'
Set-CellText B437 'G03'
Set-CellText C437 'Holdings Exit Automation MVP: Engine + order creation + reconciliation'
Set-CellText D437 'S34_G03_TB008'
Set-CellText E437 'Backend/Engine: Reconciliation job for ORDER_CREATED subscriptions: if pending order becomes EXECUTED -> COMPLETED; if CANCELLED/FAILED/REJECTED_* -> PAUSED with last_error.'
Set-CellText G437 'planned'

# Row 438
Set-CellText A438 'S34'
Set-CellText B438 'G04'
Set-CellText C438 'Holdings Exit Automation MVP: Frontend UX (Holdings + Managed Exits)'
Set-CellText D438 'S34_G04_TF001'
Set-CellText E438 'Frontend: Add holdings exit subscriptions service client (CRUD + events) with typed schemas.'
Set-CellText G438 'planned'

# Row 439
Set-CellText A439 'S34'
Set-CellText B439 'G04'
Set-CellText C439 'Holdings Exit Automation MVP: Frontend UX (Holdings + Managed Exits)'
Set-CellText D439 'S34_G04_TF002'
Set-CellText E439 'Frontend/Holdings: Add "Exit Plan" action per holding row and show subscription badge (ACTIVE/PAUSED/ORDER_CREATED/ERROR).'
Set-CellText G439 'planned'

# Row 440
Set-CellText A440 'S34'
Set-CellText B440 'G04'
Set-CellText C440 'Holdings Exit Automation MVP: Frontend UX (Holdings + Managed Exits)'
Set-CellText D440 'S34_G04_TF003'
Set-CellText E440 'Frontend/Holdings: Build Exit Plan dialog (single-leg) with size selector (% of position vs qty), trigger type selector, computed preview, and validation; default execution=MANUAL.'
Set-CellText G440 'planned'

# Row 441
Set-CellText A441 'S34'
Set-CellText B441 'G04'
Set-CellText C441 'Holdings Exit Automation MVP: Frontend UX (Holdings + Managed Exits)'
Set-CellText D441 'S34_G04_TF004'
Set-CellText E441 'Frontend/Managed Exits: Add a "Holdings exits" tab/panel (reuse ManagedRiskPage pattern) listing subscriptions + actions (pause/resume/delete) + view events.'
Set-CellText G441 'planned'

# Row 442
Set-CellText A442 'S34'
Set-CellText B442 'G04'
Set-CellText C442 'Holdings Exit Automation MVP: Frontend UX (Holdings + Managed Exits)'
Set-CellText D442 'S34_G04_TF005'
Set-CellText E442 'Frontend/Queue: Improve labeling/sorting for subscription-created orders (client_order_id starts with HEX:) so subscription exits are easy to spot (and can be shown first when needed).'
Set-CellText G442 'planned'

# Row 443
Set-CellText A443 'S34'
Set-CellText B443 'G05'
Set-CellText C443 'Holdings Exit Automation MVP: Tests + QA + docs'
Set-CellText D443 'S34_G05_TB001'
Set-CellText E443 'Backend tests: Unit tests for trigger evaluation + qty resolution + next_eval_at adaptive schedule (edge cases: small qty, rounding, percent).'
Set-CellText G443 'planned'

# Row 444
Set-CellText A444 'S34'
Set-CellText B444 'G05'
Set-CellText C444 'Holdings Exit Automation MVP: Tests + QA + docs'
Set-CellText D444 'S34_G05_TB002'
Set-CellText E444 'Backend tests: Integration tests for engine trigger->order creation and idempotency across reruns; monkeypatch broker holdings + quote fetch.'
Set-CellText G444 'planned'

# Row 445
Set-CellText A445 'S34'
Set-CellText B445 'G05'
Set-CellText C445 'Holdings Exit Automation MVP: Tests + QA + docs'
Set-CellText D445 'S34_G05_TB003'
Set-CellText E445 'Backend tests: Conflict test - subscription trigger + TradingView SELL on same symbol produce independent intents (both WAITING), with subscription one clearly labeled and no auto-dispatch.'
Set-CellText G445 'planned'

# Row 446
Set-CellText A446 'S34'
Set-CellText B446 'G05'
Set-CellText C446 'Holdings Exit Automation MVP: Tests + QA + docs'
Set-CellText D446 'S34_G05_TF001'
Set-CellText E446 'Frontend tests: Exit Plan dialog validation + create flow; Managed Exits list actions (pause/resume/delete).'
Set-CellText G446 'planned'

# Row 447
Set-CellText A447 'S34'
Set-CellText B447 'G05'
Set-CellText C447 'Holdings Exit Automation MVP: Tests + QA + docs'
Set-CellText D447 'S34_G05_TD001'
Set-CellText E447 'Docs/QA: Add QA checklist + scenario matrix for holdings exit automation (trigger, cancel, reject, conflict with TV, broker outage).'
Set-CellText G447 'planned'

# Row 448
Set-CellText A448 'S34'
Set-CellText B448 'G05'
Set-CellText C448 'Holdings Exit Automation MVP: Tests + QA + docs'
Set-CellText D448 'S34_G05_TD002'
Set-CellText E448 'Docs/Ops: Add rollout + config notes (feature flags, recommended poll interval, safe defaults, Postgres concurrency notes).'
Set-CellText G448 'planned'

# Row 449
Set-CellText A449 'S34'
Set-CellText B449 'G05'
Set-CellText C449 'Holdings Exit Automation MVP: Tests + QA + docs'
Set-CellText D449 'S34_G05_TD003'
Set-CellText E449 'Repo hygiene: Add .gitignore rule for `*.db-journal` and other SQLite transient files to prevent accidental commits.'
Set-CellText G449 'planned'

